# Rewrite speaker notes to concise reference style across 8 slides.
# Each notes page has a "Notes Placeholder" shape (index 2) whose
# TextFrame.TextRange.Text holds the single run of speaker-note text.

$p = $ppt.ActivePresentation

function Set-NotesText($SlideIndex, $NewText) {
    $slide = $p.Slides.Item($SlideIndex)
    $notesPage = $slide.NotesPage
    for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
        $shape = $notesPage.Shapes.Item($i)
        if ($shape.Name -eq "Notes Placeholder 2") {
            $shape.TextFrame.TextRange.Text = $NewText
            break
        }
    }
}

Set-NotesText 1 "Session 1 of AI Academy."

Set-NotesText 10 "Specificity in = specificity out. No new info — just described actual requirements."

Set-NotesText 14 "Keywords significantly change AI behavior. Leading questions → sycophantic agreement. Permit uncertainty."

Set-NotesText 16 "Specify format, tone, length, audience. Same principle as briefing a designer."

Set-NotesText 18 "Emotional peak. Full minute to compare. Transformation is visceral — self-driven."

Set-NotesText 21 "Cliffhanger from main deck. All iterations changed words, not context. Next session fixes this."

Set-NotesText 3 "Quick framing. Audience = Explorers. Used ChatGPT, vague prompts, mediocre results."

Set-NotesText 4 "Minimal level set from main deck. Don’t dwell. ~30 seconds."
